$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 6734
$ws.Cells.Item(3, 6).Value = 806
$ws.Cells.Item(5, 6).Value = 136
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(7, 6).Value = 713
$ws.Cells.Item(8, 6).Value = 713
$ws.Cells.Item(9, 6).Value = 10
$ws.Cells.Item(10, 6).Value = 222
$ws.Cells.Item(11, 6).Value = 19
$ws.Cells.Item(12, 6).Value = 1096
$ws.Cells.Item(13, 6).Value = 852
$ws.Cells.Item(14, 6).Value = 698
$ws.Cells.Item(15, 6).Value = 8
$ws.Cells.Item(16, 6).Value = 1012
$ws.Cells.Item(17, 6).Value = 1332
$ws.Cells.Item(18, 6).Value = 45
$ws.Cells.Item(20, 6).Value = 530
$ws.Cells.Item(21, 6).Value = 550
$ws.Cells.Item(24, 6).Value = 364
$ws.Cells.Item(26, 6).Value = 1488
$ws.Cells.Item(27, 6).Value = 723
$ws.Cells.Item(28, 6).Value = 516
$ws.Cells.Item(30, 6).Value = 448
$ws.Cells.Item(32, 6).Value = 7
$ws.Cells.Item(34, 6).Value = 257
$ws.Cells.Item(35, 6).Value = 2350
$ws.Cells.Item(37, 6).Value = 1224
$ws.Cells.Item(38, 6).Value = 439
$ws.Cells.Item(40, 6).Value = 3847

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 748
$ws.Cells.Item(12, 6).Value = 641
$ws.Cells.Item(18, 6).Value = 335
$ws.Cells.Item(19, 6).Value = 4125
$ws.Cells.Item(25, 6).Value = 230
$ws.Cells.Item(26, 6).Value = 247
$ws.Cells.Item(27, 6).Value = 107
$ws.Cells.Item(29, 6).Value = 225

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 49
$ws.Cells.Item(5, 6).Value = 1629
$ws.Cells.Item(8, 6).Value = 952

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 1629
$ws.Cells.Item(7, 6).Value = 952
$ws.Cells.Item(9, 6).Value = 6734
$ws.Cells.Item(11, 6).Value = 806
$ws.Cells.Item(12, 6).Value = 748
$ws.Cells.Item(13, 6).Value = 136
$ws.Cells.Item(14, 6).Value = 713
$ws.Cells.Item(15, 6).Value = 713
$ws.Cells.Item(16, 6).Value = 222
$ws.Cells.Item(17, 6).Value = 19
$ws.Cells.Item(18, 6).Value = 1096
$ws.Cells.Item(19, 6).Value = 852
$ws.Cells.Item(21, 6).Value = 698
$ws.Cells.Item(25, 6).Value = 1012
$ws.Cells.Item(26, 6).Value = 1332
$ws.Cells.Item(27, 6).Value = 45
$ws.Cells.Item(29, 6).Value = 530
$ws.Cells.Item(30, 6).Value = 550
$ws.Cells.Item(32, 6).Value = 335
$ws.Cells.Item(33, 6).Value = 364
$ws.Cells.Item(35, 6).Value = 1488
$ws.Cells.Item(37, 6).Value = 723
$ws.Cells.Item(38, 6).Value = 516
$ws.Cells.Item(40, 6).Value = 448
$ws.Cells.Item(44, 6).Value = 257
$ws.Cells.Item(45, 6).Value = 2350
$ws.Cells.Item(46, 6).Value = 225
$ws.Cells.Item(49, 6).Value = 1224
$ws.Cells.Item(50, 6).Value = 439
$ws.Cells.Item(51, 6).Value = 3847
